# Update "想去人数" (F column) figures across the four worksheets to match
# the freshly generated gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 196
$ws1.Range("F15").Value = 529
$ws1.Range("F16").Value = 924
$ws1.Range("F17").Value = 85527
$ws1.Range("F18").Value = 85527
$ws1.Range("F22").Value = 39357
$ws1.Range("F29").Value = 36
$ws1.Range("F32").Value = 769
$ws1.Range("F36").Value = 5610
$ws1.Range("F39").Value = 12

# --- 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F45").Value = 357

# --- 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 740
$ws3.Range("F5").Value = 610
$ws3.Range("F6").Value = 657
$ws3.Range("F8").Value = 117

# --- 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 610
$ws4.Range("F6").Value = 657
$ws4.Range("F7").Value = 657
$ws4.Range("F14").Value = 196
$ws4.Range("F17").Value = 117
$ws4.Range("F26").Value = 529
$ws4.Range("F27").Value = 85527
$ws4.Range("F30").Value = 39357
$ws4.Range("F40").Value = 5610
